$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-12 -> 2023-09-13) for every data row (rows 2 through 46).
$range = $ws.Range("C2:C46")
$range.Value = 45182
